$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 187 (pushes the existing rows 187-306 down to 188-307,
# growing the used range from A1:R306 to A1:R307).
$ws.Rows.Item(187).Insert()

# Populate the newly-inserted row with the new price-report record.
$ws.Range("A187").Value = 4
$ws.Range("B187").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C187").Value = "Los Lagos"
$ws.Range("D187").Value = 44879
$ws.Range("E187").Value = 10
$ws.Range("F187").Value = 100112044
$ws.Range("G187").Value = "Perejil"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 70
$ws.Range("K187").Value = 6000
$ws.Range("L187").Value = 6000
$ws.Range("M187").Value = 6000
$ws.Range("N187").Value = "$/docena de atados (2 kilos)"
$ws.Range("O187").Value = "Región de La Araucanía"
$ws.Range("P187").Value = 3000
$ws.Range("Q187").Value = 2
$ws.Range("R187").Value = "Hortaliza"
